$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.004.45"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.89"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.88"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.47"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.21"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.37"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.009.73"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.44"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.22"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.52"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.37"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123.18"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  -5.55%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0154"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.30"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.764.68"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.756"
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("E44").Value = "  -5.80%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0523"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.47"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  +0.51%  "
